$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 39: Skander Hathroubi
$ws.Cells.Item(39, 1).Value = "Skander"
$ws.Cells.Item(39, 2).Value = "Hathroubi"
$ws.Cells.Item(39, 3).Value = "Université de Strasbourg"
$ws.Cells.Item(39, 4).Value = "France"
$ws.Cells.Item(39, 5).Value = "hMHGOioAAAAJ"
$ws.Cells.Item(39, 6).Value = "M"
$ws.Cells.Item(39, 7).Value = 1985
$ws.Cells.Item(39, 8).Value = "Médecine, Biologie et Sciences de la Santé"

# New row 40: Omayma Missawi
$ws.Cells.Item(40, 1).Value = "Omayma"
$ws.Cells.Item(40, 2).Value = "Missawi"
$ws.Cells.Item(40, 3).Value = "Université de Namur"
$ws.Cells.Item(40, 4).Value = "Belgique"
$ws.Cells.Item(40, 5).Value = "qMrkBQsAAAAJ"
$ws.Cells.Item(40, 6).Value = "F"
$ws.Cells.Item(40, 7).Value = 1992
$ws.Cells.Item(40, 8).Value = "Médecine, Biologie et Sciences de la Santé"

# Match the style used in column F (Genre) for the rest of the data rows, by
# copying the formatting from the cell directly above instead of recreating it,
# so no new/duplicate style entries are introduced.
$ws.Range("F38").Copy()
$ws.Range("F39:F40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update view state to reflect scrolled position / active cell after edit
$ws.Application.Goto($ws.Range("A22"))
$ws.Range("H41").Select()
